# Generate Report for Handback
#
# The localization-status report previously only recorded handoff info
# (Source File, Status, Latest Handoff File/Datetime) and left the
# handback columns (Latest Target File, Latest Handback File,
# Latest Handback DateTime) empty/default. Now that the handback for the
# two tracked files has completed, this script stamps the status as
# "Handed back: in sync with en-US", fills in the Latest Target File /
# Latest Handback File columns (pointing at the same source/handoff
# files, now confirmed back in sync) and records the real handback
# datetime in place of the "0001-01-01 00:00:00" placeholder.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: just the status text changes ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

function Update-LangSheet {
    param($sheet, $md1Url, $xlf2Url, $md3Url, $xlf3Url, $handbackTime)

    # Row 2 - b1beacd1-... entry
    $sheet.Range("B2").Value = $newStatus
    $srcA2 = $sheet.Range("A2").Text
    $srcC2 = $sheet.Range("C2").Text
    $sheet.Hyperlinks.Add($sheet.Range("E2"), $md1Url, "", "", $srcA2) | Out-Null
    $sheet.Hyperlinks.Add($sheet.Range("F2"), $xlf2Url, "", "", $srcC2) | Out-Null
    $sheet.Range("G2").Value = $handbackTime

    # Row 3 - f57829ef-... entry
    $sheet.Range("B3").Value = $newStatus
    $srcA3 = $sheet.Range("A3").Text
    $srcC3 = $sheet.Range("C3").Text
    $sheet.Hyperlinks.Add($sheet.Range("E3"), $md3Url, "", "", $srcA3) | Out-Null
    $sheet.Hyperlinks.Add($sheet.Range("F3"), $xlf3Url, "", "", $srcC3) | Out-Null
    $sheet.Range("G3").Value = $handbackTime
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Update-LangSheet `
    $zhcn `
    "https://github.com/OpenLocalizationTest/oltest/blob/83e2019e13564f268545e4158350c5af2b304bd6/e2e/b1beacd1-6e4f-45bc-8352-436126d411ef.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cc82109d07e207793624645ebbbe9ff2bb2fb1e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b1beacd1-6e4f-45bc-8352-436126d411ef.786bfc1b6fdee835f5a90e03138bb0dbc4f3f712.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/83e2019e13564f268545e4158350c5af2b304bd6/e2e/f57829ef-8225-4cd2-890a-00cf54974452.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cc82109d07e207793624645ebbbe9ff2bb2fb1e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f57829ef-8225-4cd2-890a-00cf54974452.5029642f4f9b71fb3c5ba5befa9b2878be2109e9.zh-cn.xlf" `
    "2016-03-08 06:35:03"

$dede = $wb.Worksheets.Item("de-de")
Update-LangSheet `
    $dede `
    "https://github.com/OpenLocalizationTest/oltest/blob/83e2019e13564f268545e4158350c5af2b304bd6/e2e/b1beacd1-6e4f-45bc-8352-436126d411ef.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fefdaf0de0c35a8e172d42ff59f7c21549db0ccb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b1beacd1-6e4f-45bc-8352-436126d411ef.786bfc1b6fdee835f5a90e03138bb0dbc4f3f712.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/83e2019e13564f268545e4158350c5af2b304bd6/e2e/f57829ef-8225-4cd2-890a-00cf54974452.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fefdaf0de0c35a8e172d42ff59f7c21549db0ccb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f57829ef-8225-4cd2-890a-00cf54974452.5029642f4f9b71fb3c5ba5befa9b2878be2109e9.de-de.xlf" `
    "2016-03-08 06:35:21"
